$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.8180139725122046
$ws.Cells.Item(2, 3).Value = -0.4385921514182543
$ws.Cells.Item(2, 4).Value = 0.3092605588155036
$ws.Cells.Item(2, 5).Value = 0.3204530529485543
$ws.Cells.Item(3, 2).Value = 0.05413808103344353
$ws.Cells.Item(3, 3).Value = 0.6752354387775248
$ws.Cells.Item(3, 4).Value = 0.3671020368216743
$ws.Cells.Item(3, 5).Value = 0.3782945309547249
$ws.Cells.Item(4, 2).Value = -0.4699902090725725
$ws.Cells.Item(4, 3).Value = 0.4749393797381012
$ws.Cells.Item(4, 4).Value = -0.02772933836605468
$ws.Cells.Item(4, 5).Value = -0.01653684423300405
$ws.Cells.Item(5, 2).Value = 0.5467038660876931
$ws.Cells.Item(5, 3).Value = -0.003159581291080071
$ws.Cells.Item(5, 4).Value = 0.1117741737900311
$ws.Cells.Item(5, 5).Value = 0.1229666679230817
$ws.Cells.Item(6, 2).Value = 0.6094430382877407
$ws.Cells.Item(6, 3).Value = -0.7914211707358625
$ws.Cells.Item(6, 4).Value = -0.5088249777282632
$ws.Cells.Item(6, 5).Value = -0.4976324835952126
$ws.Cells.Item(7, 2).Value = 0.26044687922147
$ws.Cells.Item(7, 3).Value = -0.3999230529977498
$ws.Cells.Item(7, 4).Value = -0.158933041270463
$ws.Cells.Item(7, 5).Value = -0.1477405471374124
$ws.Cells.Item(8, 2).Value = 0.7398322420886418
$ws.Cells.Item(8, 3).Value = -0.1316216516568331
$ws.Cells.Item(8, 4).Value = 0.04982596342212067
$ws.Cells.Item(8, 5).Value = 0.0610184575551713
$ws.Cells.Item(9, 2).Value = 0.8974631314321124
$ws.Cells.Item(9, 3).Value = -0.537200097726501
$ws.Cells.Item(9, 4).Value = -0.3502289269760306
$ws.Cells.Item(9, 5).Value = -0.33903643284298
$ws.Cells.Item(10, 2).Value = -0.9313160630050259
$ws.Cells.Item(10, 3).Value = -0.3640289314639717
$ws.Cells.Item(10, 4).Value = 0.3318350354796737
$ws.Cells.Item(10, 5).Value = 0.3430275296127244
$ws.Cells.Item(11, 2).Value = 0.649691719193338
$ws.Cells.Item(11, 3).Value = -0.6484078981324237
$ws.Cells.Item(11, 4).Value = -0.409573208095256
$ws.Cells.Item(11, 5).Value = -0.3983807139622054
$ws.Cells.Item(12, 2).Value = 0.4629325164997189
$ws.Cells.Item(12, 3).Value = 0.2341588585026415
$ws.Cells.Item(12, 4).Value = 0.2800519377065192
$ws.Cells.Item(12, 5).Value = 0.2912444318395698
$ws.Cells.Item(13, 2).Value = -0.7802850270179154
$ws.Cells.Item(13, 3).Value = -0.6329263462578876
$ws.Cells.Item(13, 4).Value = 0.4024898801893085
$ws.Cells.Item(13, 5).Value = 0.4136823743223592
$ws.Cells.Item(14, 2).Value = -0.9702287095938849
$ws.Cells.Item(14, 3).Value = -0.6324332265968937
$ws.Cells.Item(14, 4).Value = 0.5698195966884069
$ws.Cells.Item(14, 5).Value = 0.5810120908214575
$ws.Cells.Item(15, 2).Value = -0.5174014858441027
$ws.Cells.Item(15, 3).Value = -0.6791428774561667
$ws.Cells.Item(15, 4).Value = 0.2103699414174278
$ws.Cells.Item(15, 5).Value = 0.2215624355504784
$ws.Cells.Item(16, 2).Value = 0.8651340609620961
$ws.Cells.Item(16, 3).Value = -0.5896833576772957
$ws.Cells.Item(16, 4).Value = -0.4018149287499378
$ws.Cells.Item(16, 5).Value = -0.3906224346168872
$ws.Cells.Item(17, 2).Value = -0.9731474481201052
$ws.Cells.Item(17, 3).Value = -0.1696907746282261
$ws.Cells.Item(17, 4).Value = 0.2115857261641171
$ws.Cells.Item(17, 5).Value = 0.2227782202971677
$ws.Cells.Item(18, 2).Value = -0.08188190065829004
$ws.Cells.Item(18, 3).Value = 0.4823250952846854
$ws.Cells.Item(18, 4).Value = 0.1561047619386784
$ws.Cells.Item(18, 5).Value = 0.167297256071729
$ws.Cells.Item(19, 2).Value = 0.1980134207212292
$ws.Cells.Item(19, 3).Value = -0.8159509586557261
$ws.Cells.Item(19, 4).Value = -0.2998801834921427
$ws.Cells.Item(19, 5).Value = -0.2886876893590921
$ws.Cells.Item(20, 2).Value = -0.7776166949092667
$ws.Cells.Item(20, 3).Value = 0.9256884273401516
$ws.Cells.Item(20, 4).Value = -0.1482464041885797
$ws.Cells.Item(20, 5).Value = -0.1370539100555291
$ws.Cells.Item(21, 2).Value = 0.203215898749068
$ws.Cells.Item(21, 3).Value = 0.5826619044744215
$ws.Cells.Item(21, 4).Value = 0.4094657172915553
$ws.Cells.Item(21, 5).Value = 0.420658211424606
$ws.Cells.Item(22, 2).Value = -0.4698230886947183
$ws.Cells.Item(22, 3).Value = 0.6746659055476678
$ws.Cells.Item(22, 4).Value = 0.004334450532910206
$ws.Cells.Item(22, 5).Value = 0.01552694466596083
$ws.Cells.Item(23, 2).Value = 0.7090078928075463
$ws.Cells.Item(23, 3).Value = 0.21635342230348
$ws.Cells.Item(23, 4).Value = 0.4001160235817616
$ws.Cells.Item(23, 5).Value = 0.4113085177148122
$ws.Cells.Item(24, 2).Value = 0.1960798249591653
$ws.Cells.Item(24, 3).Value = 0.04203233153523667
$ws.Cells.Item(24, 4).Value = 0.04850960301962629
$ws.Cells.Item(24, 5).Value = 0.05970209715267692
$ws.Cells.Item(25, 2).Value = 0.4402964239131835
$ws.Cells.Item(25, 3).Value = 0.727485269050417
$ws.Cells.Item(25, 4).Value = 0.7657045903388009
$ws.Cells.Item(25, 5).Value = 0.7768970844718516
